$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: UI entry for "rpg-gui-construction-kit-v10"
$ws.Hyperlinks.Add(
    $ws.Range("B2"),
    "https://opengameart.org/content/rpg-gui-construction-kit-v10",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "https://opengameart.org/content/rpg-gui-construction-kit-v10"
) | Out-Null
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("A2").Value = "UI"

# Row 3: Options entry for "0-ad-gui-elements"
$ws.Hyperlinks.Add(
    $ws.Range("B3"),
    "https://opengameart.org/content/0-ad-gui-elements",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "https://opengameart.org/content/0-ad-gui-elements"
) | Out-Null
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("A3").Value = "Options"

# Update selection to A4
$ws.Range("A4").Select() | Out-Null
